# Apply the changes described by the diff:
#  1) D655 changes from "Duesseldorf" to "Neuss"
#  2) Rows 855-866 (previously blank placeholder rows) get filled in with
#     new event data, including new Instagram hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Simple value change: D655 -> "Neuss"
# ---------------------------------------------------------------------
$ws.Range("D655").Value = "Neuss"

# ---------------------------------------------------------------------
# 2) Fill in rows 855-866 with new data.
#    Each row currently has no value (but already has the right base
#    cell styles: A=date style, B:E=blank/general style). We copy the
#    format from the previous, fully-populated row (854) down through
#    each new row so that the text columns (B:E) pick up the normal
#    "text" style, then set the actual values.
# ---------------------------------------------------------------------

$rows = @(
    @{ Row = 855; Date = 45976; Event = "CTRL";                         Location = "SNRS";        City = "Dortmund";  Link = "https://www.instagram.com/reel/DQEy_nBjBM8/?igsh=Z25xamdybWlmazln" },
    @{ Row = 856; Date = 45976; Event = "POLYAMOR";                     Location = "Junkyard";     City = "Dortmund";  Link = "https://www.instagram.com/p/DPWl3GajJn6/?igsh=MThzcTltM3B6dm0xbA==" },
    @{ Row = 857; Date = 45954; Event = "HIGH VOLTAGE";                 Location = "Elektroküche"; City = "Köln";      Link = "https://www.instagram.com/reel/DPWfdmhjZEj/?igsh=MWQ4MnB5em5wOTF5MA==" },
    @{ Row = 858; Date = 45953; Event = "POWER 4 HOUR";                 Location = "SNRS";        City = "Dortmund";  Link = "https://www.instagram.com/reel/DQDFj1IDNUi/?igsh=MXZoajUxMHNwbWF3Zg==" },
    @{ Row = 859; Date = 45970; Event = "HÖR x HAFEN7";                 Location = "hafen7";       City = "Neuss";     Link = "https://www.instagram.com/reel/DPoLBrYiLNl/?igsh=MXc4ZXFtd2kydG1kYg==" },
    @{ Row = 860; Date = 45961; Event = "GVOE HALLOWEEN";               Location = "Ground Zero";  City = "Essen";     Link = "https://www.instagram.com/p/DPeYrEWgndb/?igsh=OGlmOWFwaXRuYXcy" },
    @{ Row = 861; Date = 45952; Event = "#MITTWOCHENENDE";              Location = "Odonien";      City = "Köln";      Link = "https://www.instagram.com/p/DP4GF97jLY7/?igsh=bGsyZXhtNWp2MXQw" },
    @{ Row = 862; Date = 45961; Event = "HALLOWEEN SPECIAL FREE RAVE";  Location = "Zimmermanns";  City = "Köln";      Link = "https://www.instagram.com/reel/DP6lXMFiIEd/?igsh=ZzlrYTlheDF0MzR4" },
    @{ Row = 863; Date = 45996; Event = "RAVE EXPANSION";               Location = "Dings";        City = "Köln";      Link = "https://www.instagram.com/reel/DP1Yh3ICKkw/?igsh=eTByMzVsa3RkcDc4" },
    @{ Row = 864; Date = 45961; Event = "BLACKCELL TECHNO";             Location = "Lessie";       City = "Aachen";    Link = "https://www.instagram.com/reel/DPbHEMrDFM0/?igsh=OHZwZ2tpamJzeG8z" },
    @{ Row = 865; Date = 45968; Event = "RAVE IN ODONIEN";              Location = "Odonien";      City = "Köln";      Link = "https://www.instagram.com/reel/DP_NwogiDPX/?igsh=MTRjc2swZDNseWJ4eg==" },
    @{ Row = 866; Date = 46004; Event = "CLUB NIGHT";                   Location = "SNRS";        City = "Dortmund";  Link = "https://www.instagram.com/reel/DQErqnIDFqv/?igsh=aXdkcW45aWd3bnk2" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Copy the formatting from row 854 (the last fully populated row)
    # down onto this row so the cells use the normal text/date styles
    # instead of the empty-placeholder style.
    $ws.Range("A854:E854").Copy()
    $ws.Range("A" + $rowNum + ":E" + $rowNum).PasteSpecial(-4122)

    $ws.Range("A" + $rowNum).Value = $r.Date
    $ws.Range("B" + $rowNum).Value = $r.Event
    $ws.Range("C" + $rowNum).Value = $r.Location
    $ws.Range("D" + $rowNum).Value = $r.City

    $eCell = $ws.Range("E" + $rowNum)
    $link = $r.Link
    $eCell.Value = $link

    # Register the real hyperlink relationship (target + display text).
    $ws.Hyperlinks.Add($eCell, $link, "", "", $link)

    # Style the link text itself (underline, blue) same as the rest of
    # the existing hyperlink cells in column E. The run is intentionally
    # split into two Characters() calls (instead of one call spanning
    # the whole string) so the engine keeps the formatting as an
    # explicit rich-text run on the shared string rather than collapsing
    # it back down to a plain string with only a cell-level style.
    $len = $link.Length
    $eCell.Characters(1, $len - 1).Font.Underline = 2
    $eCell.Characters(1, $len - 1).Font.ColorIndex = 11
    $eCell.Characters($len, 1).Font.Underline = 2
    $eCell.Characters($len, 1).Font.ColorIndex = 11

    # Re-apply the plain text cell format (border/fill/number format)
    # on top, since adding the hyperlink can reset it.
    $ws.Range("E854").Copy()
    $eCell.PasteSpecial(-4122)
}
